# Update the NATMI LR-pair sheet (Col18a1-Gpc4) with recomputed TPM-derived
# statistics. Columns A-D (cluster/gene labels) and K/L are unchanged; only
# the numeric expression / specificity / edge-weight columns (E-J, M-T) move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 3;                    "F2" = 1;                    "G2" = 0.5586043333333334
    "H2" = 1.675813;             "I2" = 0.01643366487114074;  "J2" = 0.01643366487114074
    "M2" = 5.333065333333334;    "N2" = 15.999196;            "O2" = 0.1422335064894009
    "P2" = 0.1422335064894009;   "Q2" = 2.979073405149778;    "R2" = 26.811660646348
    "S2" = 0.002337417779094036; "T2" = 0.002337417779094036

    "E3" = 3;                    "F3" = 1;                    "G3" = 0.5586043333333334
    "H3" = 1.675813;             "I3" = 0.01643366487114074;  "J3" = 0.01643366487114074
    "M3" = 20.88867166666667;    "N3" = 62.666015;            "O3" = 0.5571034351455781
    "P3" = 0.5571034351455783;   "Q3" = 11.66850251057722;    "R3" = 105.016522595195
    "S3" = 0.009155251151743721; "T3" = 0.009155251151743721

    "E4" = 3;                    "F4" = 1;                    "G4" = 0.5586043333333334
    "H4" = 1.675813;             "I4" = 0.01643366487114074;  "J4" = 0.01643366487114074
    "M4" = 11.27340366666667;    "N4" = 33.820211;            "O4" = 0.3006630583650208
    "P4" = 0.3006630583650208;   "Q4" = 6.29737213961589;     "R4" = 56.676349256543
    "S4" = 0.004940995940302981; "T4" = 0.00494099594030298

    "H5" = 63.825936;            "I5" = 0.6259016025719319;   "J5" = 0.6259016025719319
    "M5" = 5.333065333333334;    "N5" = 15.999196;            "O5" = 0.1422335064894009
    "P5" = 0.1422335064894009;   "Q5" = 113.4626288830507;    "R5" = 1021.163659947456
    "S5" = 0.08902417965114132;  "T5" = 0.08902417965114132

    "H6" = 63.825936;            "I6" = 0.6259016025719319;   "J6" = 0.6259016025719319
    "M6" = 20.88867166666667;    "N6" = 62.666015;            "O6" = 0.5571034351455781
    "P6" = 0.5571034351455783;   "S6" = 0.3486919328559457;   "T6" = 0.3486919328559458

    "H7" = 63.825936;            "I7" = 0.6259016025719319;   "J7" = 0.6259016025719319
    "M7" = 11.27340366666667;    "N7" = 33.820211;            "O7" = 0.3006630583650208
    "P7" = 0.3006630583650208;   "Q7" = 239.8451803102773;    "R7" = 2158.606622792496
    "S7" = 0.1881854900648448;   "T7" = 0.1881854900648448

    "G8" = 12.157548;            "H8" = 36.472644;            "I8" = 0.3576647325569273
    "J8" = 0.3576647325569273;   "M8" = 5.333065333333334;    "N8" = 15.999196
    "O8" = 0.1422335064894009;   "P8" = 0.1422335064894009;   "Q8" = 64.83699777713601
    "R8" = 583.5329799942241;    "S8" = 0.05087190905916557;  "T8" = 0.05087190905916557

    "G9" = 12.157548;            "H9" = 36.472644;            "I9" = 0.3576647325569273
    "J9" = 0.3576647325569273;   "M9" = 20.88867166666667;    "N9" = 62.666015
    "O9" = 0.5571034351455781;   "P9" = 0.5571034351455783;   "Q9" = 253.95502844374
    "R9" = 2285.59525599366;     "S9" = 0.1992562511378887;   "T9" = 0.1992562511378888

    "G10" = 12.157548;           "H10" = 36.472644;           "I10" = 0.3576647325569273
    "J10" = 0.3576647325569273;  "M10" = 11.27340366666667;   "N10" = 33.820211
    "O10" = 0.3006630583650208;  "P10" = 0.3006630583650208;  "Q10" = 137.056946200876
    "R10" = 1233.512515807884;   "S10" = 0.107536572359873;   "T10" = 0.107536572359873
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
